$d = $word.ActiveDocument

# Mapping of old "a÷b=" expressions to new ones, in document order.
$replacements = @(
    @("88÷5=", "55÷7="),
    @("48÷5=", "23÷8="),
    @("54÷7=", "73÷4="),
    @("44÷6=", "10÷7="),
    @("36÷8=", "89÷4="),
    @("31÷6=", "21÷9="),
    @("58÷7=", "75÷5="),
    @("43÷9=", "77÷5="),
    @("82÷3=", "23÷4="),
    @("69÷3=", "73÷8="),
    @("39÷7=", "47÷8="),
    @("98÷6=", "27÷9="),
    @("68÷6=", "14÷3="),
    @("69÷9=", "60÷2="),
    @("37÷3=", "24÷9="),
    @("19÷2=", "12÷4="),
    @("85÷7=", "46÷4="),
    @("24÷8=", "94÷6="),
    @("42÷8=", "93÷3="),
    @("35÷2=", "17÷7="),
    @("10÷4=", "64÷2="),
    @("99÷9=", "64÷2="),
    @("26÷4=", "76÷3="),
    @("22÷8=", "51÷2="),
    @("76÷8=", "14÷3=")
)

foreach ($pair in $replacements) {
    $oldText = $pair[0]
    $newText = $pair[1]

    $range = $d.Content
    $range.Find.Execute($oldText, $true, $false, $false, $false, $false, `
                         $true, 1, $false, $newText, 2)
}
